$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "Mint chocolate chip"
$ws.Range("F6").Value = "onion"

$ws.Range("F6").Select()
